# Update the handoff/handback timestamps for the
# "50e76ac2-95c0-4c90-8b1e-1b9e89a85865" file rows after re-generating the
# handback status report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 3 & 4 both reference the 50e76ac2... handoff file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-23 06:23:10"
$wsZhCn.Range("E4").Value = "2016-03-23 06:23:10"
$wsZhCn.Range("H3").Value = "2016-03-23 06:23:34"
$wsZhCn.Range("H4").Value = "2016-03-23 06:23:34"

# de-de sheet: rows 3 & 4 both reference the 50e76ac2... handoff file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-23 06:23:14"
$wsDeDe.Range("E4").Value = "2016-03-23 06:23:14"
$wsDeDe.Range("H3").Value = "2016-03-23 06:23:41"
$wsDeDe.Range("H4").Value = "2016-03-23 06:23:41"
